$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row (row 11): right-answer mark value changes from 3 to 5
$ws.Range("B11").Value = 5

# "Total" row (row 12): total correct marks = (No. Right) * (Marking Right) = 22 * 5 = 110
$ws.Range("B12").Value = 110

# Corr/total marks label e.g. "63/84" -> "110/140" (110 correct out of 22*Max(28)=140)
$ws.Range("E12").Value = "110/140"
